# Apply fixes to the dictionary template workbook:
# 1. Remove the leftover "testDataset / var1" test row from the Variables sheet.
# 2. Rename the "name" header on the "Variable values" sheet to "variable".

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Rows.Item(15).Delete()

$wsValues = $wb.Worksheets.Item("Variable values")
$wsValues.Range("C1").Value = "variable"
